$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "M.Tech (Artificial Intelligence & Machine Learning)" section header row (110)
# and its repeated table-header row (111) first (higher row numbers first so earlier
# row numbers stay stable while deleting).
$ws.Rows("110:111").Delete()

# Delete the "Computer Science & Business System (CSBS)" section header row (105)
# and its repeated table-header row (106).
$ws.Rows("105:106").Delete()

# Renumber the S.No. column for the rows that used to restart counting at 1 in each
# sub-section, so the whole table is now one continuously numbered list.
$ws.Range("A105").Value = 104
$ws.Range("A106").Value = 105
$ws.Range("A107").Value = 106
$ws.Range("A108").Value = 107
$ws.Range("A109").Value = 108
$ws.Range("A110").Value = 109
